# Add two new columns, I0 and IF, to the right of the existing IP (H) column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting already used by column H (bold, bordered,
# centered / top-aligned) for the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows --------------------------------------------------------------
# I0 is 1 for every game, and IF mirrors the existing IP (column H) value --
# except for row 25, which is a special case (I0=4, IF=8).
foreach ($r in 2..27) {
    $ip = $ws.Cells.Item($r, 8).Value()

    if ($r -eq 25) {
        $i0 = 4
        $if = 8
    } else {
        $i0 = 1
        $if = $ip
    }

    $ws.Cells.Item($r, 9).Value = $i0
    $ws.Cells.Item($r, 10).Value = $if
}
